# Scheduled-runner update: refresh Leve profit-tracking values
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1068.1489
$ws.Range("I15").Value = 1068.1489
$ws.Range("K15").Value = 3204.4467
$ws.Range("M15").Value = -3035.4467

# row 33
$ws.Range("H33").Value = 173.53847
$ws.Range("I33").Value = 155.6
$ws.Range("J33").Value = 233.33333
$ws.Range("K33").Value = 155.6
$ws.Range("L33").Value = 233.33333
$ws.Range("M33").Value = 73.40000000000001
$ws.Range("N33").Value = -691.3333299999999

# row 86
$ws.Range("H86").Value = 32909.355
$ws.Range("I86").Value = 2773.1
$ws.Range("J86").Value = 108250
$ws.Range("K86").Value = 2773.1
$ws.Range("L86").Value = 108250
$ws.Range("M86").Value = -1650.1
$ws.Range("N86").Value = -110496

# row 89
$ws.Range("H89").Value = 32909.355
$ws.Range("I89").Value = 2773.1
$ws.Range("J89").Value = 108250
$ws.Range("K89").Value = 13865.5
$ws.Range("L89").Value = 541250
$ws.Range("M89").Value = -8249.5
$ws.Range("N89").Value = -552482

# row 96
$ws.Range("H96").Value = 1423.875
$ws.Range("I96").Value = 1198.75
$ws.Range("J96").Value = 1649
$ws.Range("K96").Value = 3596.25
$ws.Range("L96").Value = 4947
$ws.Range("M96").Value = -2223.25
$ws.Range("N96").Value = -7693

# row 103
$ws.Range("H103").Value = 63879.145
$ws.Range("I103").Value = 200500
$ws.Range("J103").Value = 9230.799999999999
$ws.Range("K103").Value = 601500
$ws.Range("L103").Value = 27692.4
$ws.Range("M103").Value = -600914
$ws.Range("N103").Value = -28864.4

# row 112
$ws.Range("H112").Value = 18753.334
$ws.Range("I112").Value = 3085
$ws.Range("J112").Value = 20612.29
$ws.Range("K112").Value = 9255
$ws.Range("L112").Value = 61836.87
$ws.Range("M112").Value = -8147
$ws.Range("N112").Value = -64052.87

# row 113
$ws.Range("H113").Value = 20854.727

# row 137
$ws.Range("H137").Value = 1043386.5
$ws.Range("I137").Value = 3436959.2
$ws.Range("J137").Value = 17569.572
$ws.Range("K137").Value = 10310877.6
$ws.Range("L137").Value = 52708.716
$ws.Range("M137").Value = -10308327.6
$ws.Range("N137").Value = -57808.716

# row 138
$ws.Range("H138").Value = 9092.645
$ws.Range("I138").Value = 13714.143
$ws.Range("J138").Value = 8241.315000000001
$ws.Range("K138").Value = 41142.429
$ws.Range("L138").Value = 24723.945
$ws.Range("M138").Value = -36002.429
$ws.Range("N138").Value = -35003.945

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 1686
$ws.Range("I32").Value = 1724.1791
$ws.Range("K32").Value = 1724.1791
$ws.Range("M32").Value = -1437.1791

# row 36
$ws.Range("H36").Value = 19506.5
$ws.Range("I36").Value = 16513
$ws.Range("K36").Value = 16513
$ws.Range("M36").Value = -16167

# row 61
$ws.Range("H61").Value = 2798
$ws.Range("I61").Value = 1153
$ws.Range("J61").Value = 6499.25
$ws.Range("K61").Value = 1153
$ws.Range("L61").Value = 6499.25
$ws.Range("M61").Value = -941
$ws.Range("N61").Value = -6923.25

# row 97
$ws.Range("H97").Value = 6670725
$ws.Range("I97").Value = 4566.24
$ws.Range("K97").Value = 4566.24
$ws.Range("M97").Value = -4070.24

# row 110
$ws.Range("H110").Value = 1569.7059
$ws.Range("I110").Value = 726.9
$ws.Range("K110").Value = 726.9
$ws.Range("M110").Value = 1318.1

# row 122
$ws.Range("H122").Value = 247857.52
$ws.Range("I122").Value = 1655.4595
$ws.Range("J122").Value = 703331.3
$ws.Range("K122").Value = 4966.3785
$ws.Range("L122").Value = 2109993.9
$ws.Range("M122").Value = -2516.3785
$ws.Range("N122").Value = -2114893.9

# row 132
$ws.Range("H132").Value = 2320.75
$ws.Range("I132").Value = 2037.4736
$ws.Range("J132").Value = 4114.8335
$ws.Range("K132").Value = 6112.4208
$ws.Range("L132").Value = 12344.5005
$ws.Range("M132").Value = -3582.4208
$ws.Range("N132").Value = -17404.5005

# row 136
$ws.Range("H136").Value = 2798
$ws.Range("I136").Value = 1153
$ws.Range("J136").Value = 6499.25
$ws.Range("K136").Value = 3459
$ws.Range("L136").Value = 19497.75
$ws.Range("M136").Value = -909
$ws.Range("N136").Value = -24597.75

$ws = $wb.Worksheets.Item("BSM")
# row 64
$ws.Range("H64").Value = 7900
$ws.Range("I64").Value = 17279.25
$ws.Range("J64").Value = 4148.3
$ws.Range("K64").Value = 17279.25
$ws.Range("L64").Value = 4148.3
$ws.Range("M64").Value = -17054.25
$ws.Range("N64").Value = -4598.3

# row 67
$ws.Range("H67").Value = 7900
$ws.Range("I67").Value = 17279.25
$ws.Range("J67").Value = 4148.3
$ws.Range("K67").Value = 17279.25
$ws.Range("L67").Value = 4148.3
$ws.Range("M67").Value = -16499.25
$ws.Range("N67").Value = -5708.3

# row 86
$ws.Range("H86").Value = 21286
$ws.Range("I86").Value = 21286
$ws.Range("K86").Value = 21286
$ws.Range("M86").Value = -20163

# row 89
$ws.Range("H89").Value = 21286
$ws.Range("I89").Value = 21286
$ws.Range("K89").Value = 106430
$ws.Range("M89").Value = -100814

# row 107
$ws.Range("H107").Value = 3408.5
$ws.Range("I107").Value = 3466.3333
$ws.Range("J107").Value = 2888
$ws.Range("K107").Value = 3466.3333
$ws.Range("L107").Value = 2888
$ws.Range("M107").Value = -1546.3333
$ws.Range("N107").Value = -6728

# row 134
$ws.Range("H134").Value = 6574.16
$ws.Range("I134").Value = 6802.3477
$ws.Range("K134").Value = 20407.0431
$ws.Range("M134").Value = -17872.0431

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2431.8845
$ws.Range("I31").Value = 1540.9166
$ws.Range("J31").Value = 3195.5715
$ws.Range("K31").Value = 1540.9166
$ws.Range("L31").Value = 3195.5715
$ws.Range("M31").Value = -1245.9166
$ws.Range("N31").Value = -3785.5715

# row 34
$ws.Range("H34").Value = 2431.8845
$ws.Range("I34").Value = 1540.9166
$ws.Range("J34").Value = 3195.5715
$ws.Range("K34").Value = 1540.9166
$ws.Range("L34").Value = 3195.5715
$ws.Range("M34").Value = -1338.9166
$ws.Range("N34").Value = -3599.5715

# row 94
$ws.Range("H94").Value = 2814.5334
$ws.Range("I94").Value = 5595.5
$ws.Range("J94").Value = 1803.2727
$ws.Range("K94").Value = 5595.5
$ws.Range("L94").Value = 1803.2727
$ws.Range("M94").Value = -5144.5
$ws.Range("N94").Value = -2705.2727

# row 132
$ws.Range("H132").Value = 13634.977
$ws.Range("I132").Value = 4296.775
$ws.Range("J132").Value = 200399
$ws.Range("K132").Value = 12890.325
$ws.Range("L132").Value = 601197
$ws.Range("M132").Value = -10360.325
$ws.Range("N132").Value = -606257

$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value = 36000

# row 97
$ws.Range("H97").Value = 25300.4
$ws.Range("J97").Value = 1176.8
$ws.Range("L97").Value = 3530.4
$ws.Range("N97").Value = -4522.4

# row 129
$ws.Range("H129").Value = 41668580
$ws.Range("I129").Value = 1724
$ws.Range("J129").Value = 83335430
$ws.Range("K129").Value = 5172
$ws.Range("L129").Value = 250006290
$ws.Range("M129").Value = -172
$ws.Range("N129").Value = -250016290

# row 131
$ws.Range("H131").Value = 71435350
$ws.Range("J131").Value = 1965.125
$ws.Range("L131").Value = 5895.375
$ws.Range("N131").Value = -15975.375

# row 133
$ws.Range("H133").Value = 8127.125
$ws.Range("I133").Value = 4004
$ws.Range("J133").Value = 14999
$ws.Range("K133").Value = 12012
$ws.Range("L133").Value = 44997
$ws.Range("M133").Value = -6952
$ws.Range("N133").Value = -55117

$ws = $wb.Worksheets.Item("GSM")
# row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = 0

# row 126
$ws.Range("H126").Value = 25634.273
$ws.Range("I126").Value = 53999.668
$ws.Range("J126").Value = 14997.25
$ws.Range("K126").Value = 161999.004
$ws.Range("L126").Value = 44991.75
$ws.Range("M126").Value = -159529.004
$ws.Range("N126").Value = -49931.75

# row 132
$ws.Range("H132").Value = 2594.2173
$ws.Range("I132").Value = 2019.7441
$ws.Range("J132").Value = 10828.333
$ws.Range("K132").Value = 6059.2323
$ws.Range("L132").Value = 32484.999
$ws.Range("M132").Value = -3529.2323
$ws.Range("N132").Value = -37544.999

$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 36764.457
$ws.Range("I40").Value = 42975.39
$ws.Range("K40").Value = 42975.39
$ws.Range("M40").Value = -42839.39

# row 46
$ws.Range("H46").Value = 2017.15
$ws.Range("I46").Value = 1337.909
$ws.Range("J46").Value = 2847.3333
$ws.Range("K46").Value = 1337.909
$ws.Range("L46").Value = 2847.3333
$ws.Range("M46").Value = -1149.909
$ws.Range("N46").Value = -3223.3333

# row 61
$ws.Range("H61").Value = 3149.5
$ws.Range("I61").Value = 1824
$ws.Range("K61").Value = 1824
$ws.Range("M61").Value = -1622

# row 74
$ws.Range("H74").Value = 50197
$ws.Range("I74").Value = 50197
$ws.Range("K74").Value = 50197
$ws.Range("M74").Value = -49199

# row 77
$ws.Range("H77").Value = 50197
$ws.Range("I77").Value = 50197
$ws.Range("K77").Value = 150591
$ws.Range("M77").Value = -145599

# row 100
$ws.Range("H100").Value = 6965.3076
$ws.Range("I100").Value = 10114.429
$ws.Range("K100").Value = 10114.429
$ws.Range("M100").Value = -9573.429

# row 113
$ws.Range("H113").Value = 3149.5
$ws.Range("I113").Value = 1824
$ws.Range("K113").Value = 1824
$ws.Range("M113").Value = 346

# row 122
$ws.Range("H122").Value = 4962.1816
$ws.Range("I122").Value = 3759.75
$ws.Range("J122").Value = 5229.3887
$ws.Range("K122").Value = 11279.25
$ws.Range("L122").Value = 15688.1661
$ws.Range("M122").Value = -8829.25
$ws.Range("N122").Value = -20588.1661

# row 139
$ws.Range("H139").Value = 120999.4
$ws.Range("J139").Value = 97499.25
$ws.Range("L139").Value = 97499.25
$ws.Range("N139").Value = -107779.25

# row 141
$ws.Range("H141").Value = 44710
$ws.Range("J141").Value = 44710
$ws.Range("L141").Value = 44710
$ws.Range("N141").Value = -55070

$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 6807.4585
$ws.Range("I122").Value = 4492.2666
$ws.Range("K122").Value = 13476.7998
$ws.Range("M122").Value = -11026.7998

# row 136
$ws.Range("H136").Value = 3043.578
$ws.Range("I136").Value = 1996.6578
$ws.Range("J136").Value = 8726.857
$ws.Range("K136").Value = 5989.9734
$ws.Range("L136").Value = 26180.571
$ws.Range("M136").Value = -3439.9734
$ws.Range("N136").Value = -31280.571

# row 139
$ws.Range("H139").Value = 45000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280
